$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "进入备选升级牌区时：如果玩家未拥有《等级2》，则重抽本牌。<br>`n可以使用至多包含3张牌的堆叠。"
$ws.Range("D4").Value = "进入备选升级牌区时：如果玩家未拥有《等级3》，则重抽本牌。<br>可以使用至多包含4张牌的堆叠。"
$ws.Range("D5").Value = "进入备选升级牌区时：如果玩家未拥有《等级4》，则重抽本牌。<br>可以使用至多包含5张牌的堆叠。"

$ws.Rows("3:5").RowHeight = 171

$ws.Range("F5").Select() | Out-Null
